$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "optional" (yellow) header style from an existing cell (C15) so the
# new header cells S15:BI15 reuse the existing style index instead of creating new ones.
$ws.Range("C15").Copy()
$ws.Range("S15:BI15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$headers = @{
    "S15" = "biotic_relationship"
    "T15" = "chem_administration"
    "U15" = "ethnicity"
    "V15" = "extrachrom_elements"
    "W15" = "host_age"
    "X15" = "host_body_mass_index"
    "Y15" = "host_body_product"
    "Z15" = "host_body_temp"
    "AA15" = "host_diet"
    "AB15" = "host_disease"
    "AC15" = "host_family_relationship"
    "AD15" = "host_genotype"
    "AE15" = "host_height"
    "AF15" = "host_last_meal"
    "AG15" = "host_occupation"
    "AH15" = "host_phenotype"
    "AI15" = "host_pulse"
    "AJ15" = "host_sex"
    "AK15" = "host_subject_id"
    "AL15" = "host_tissue_sampled"
    "AM15" = "host_tot_mass"
    "AN15" = "ihmc_medication_code"
    "AO15" = "isolation_source"
    "AP15" = "medic_hist_perform"
    "AQ15" = "misc_param"
    "AR15" = "nose_mouth_teeth_throat_disord"
    "AS15" = "organism_count"
    "AT15" = "oxy_stat_samp"
    "AU15" = "perturbation"
    "AV15" = "rel_to_oxygen"
    "AW15" = "samp_collect_device"
    "AX15" = "samp_mat_process"
    "AY15" = "samp_salinity"
    "AZ15" = "samp_size"
    "BA15" = "samp_store_dur"
    "BB15" = "samp_store_loc"
    "BC15" = "samp_store_temp"
    "BD15" = "samp_vol_we_dna_ext"
    "BE15" = "source_material_id"
    "BF15" = "subspecf_gen_lin"
    "BG15" = "temperature"
    "BH15" = "time_last_toothbrush"
    "BI15" = "trophic_level"
}

$comments = @{
    "S15" = "Free-living or from host (define relationship)"
    "T15" = "list of chemical compounds administered to the host or site where sampling occurred, and when (e.g. antibiotics, N fertilizer, air filter); can include multiple compounds. For Chemical Entities of Biological Interest ontology (CHEBI) (v1.72), please see http://bioportal.bioontology.org/visualize/44603"
    "U15" = "ethnicity of the subject"
    "V15" = "Plasmids that have significance phenotypic consequence"
    "W15" = "Age of host at the time of sampling"
    "X15" = "body mass index of the host, calculated as weight/(height)squared"
    "Y15" = "substance produced by the host, e.g. stool, mucus, where the sample was obtained from"
    "Z15" = "core body temperature of the host when sample was collected"
    "AA15" = "type of diet depending on the sample for animals omnivore, herbivore etc., for humans high-fat, meditteranean etc.; can include multiple diet types"
    "AB15" = "Name of relevant disease, e.g. Salmonella gastroenteritis. For the controlled vocabulary, please see Human Disease Ontology, http://bioportal.bioontology.org/ontologies/1009 or MeSH, http://www.ncbi.nlm.nih.gov/mesh"
    "AE15" = "the height of subject"
    "AF15" = "content of last meal and time since feeding; can include multiple values"
    "AG15" = "most frequent job performed by subject"
    "AI15" = "resting pulse of the host, measured as beats per minute"
    "AJ15" = "Gender or physical sex of the host"
    "AK15" = "a unique identifier by which each subject can be referred to, de-identified, e.g. #131"
    "AL15" = "Type of tissue the initial sample was taken from. Controlled vocabulary, http://bioportal.bioontology.org/ontologies/1005"
    "AM15" = "total mass of the host at collection, the unit depends on host"
    "AN15" = "can include multiple medication codes"
    "AO15" = "Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived."
    "AP15" = "whether full medical history was collected"
    "AQ15" = "any other measurement performed or parameter collected, that is not listed here"
    "AR15" = "history of nose/mouth/teeth/throat disorders; can include multiple disorders"
    "AS15" = "total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts"
    "AT15" = "oxygenation status of sample"
    "AU15" = "type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types"
    "AV15" = "Aerobic or anaerobic"
    "AW15" = "Method or device employed for collecting sample"
    "AX15" = "Processing applied to the sample during or after isolation"
    "AY15" = "salinity of sample, i.e. measure of total salt concentration"
    "AZ15" = "Amount or size of sample (volume, mass or area) that was collected"
    "BA15" = "duration for which sample was stored"
    "BB15" = "location at which sample was stored, usually name of a specific freezer/room"
    "BC15" = "temperature at which sample was stored, e.g. -80"
    "BD15" = "volume (mL) or weight (g) of sample processed for DNA extraction"
    "BE15" = "unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples."
    "BF15" = "Information about the genetic distinctness of the lineage (eg., biovar, serovar)"
    "BG15" = "temperature of the sample at time of sampling"
    "BH15" = "specification of the time since last toothbrushing"
    "BI15" = "Feeding position in food chain (eg., chemolithotroph)"
}

$refs = @("S15", "T15", "U15", "V15", "W15", "X15", "Y15", "Z15", "AA15", "AB15", "AC15", "AD15", "AE15", "AF15", "AG15", "AH15", "AI15", "AJ15", "AK15", "AL15", "AM15", "AN15", "AO15", "AP15", "AQ15", "AR15", "AS15", "AT15", "AU15", "AV15", "AW15", "AX15", "AY15", "AZ15", "BA15", "BB15", "BC15", "BD15", "BE15", "BF15", "BG15", "BH15", "BI15")

foreach ($ref in $refs) {
    $ws.Range($ref).Value = $headers[$ref]
    if ($comments.ContainsKey($ref)) {
        $ws.Range($ref).AddComment($comments[$ref])
    }
}
